$d = $word.ActiveDocument

# --- Step 1: append the "Jamaican Restaurant..." sentence onto the end of
# the paragraph that currently ends "...Japan Air, Italy Hotel, and "
# (right before the _GoBack bookmark). ---
$d.Content.Find.Execute(
    "Our partnerships with businesses such as Japan Air, Italy Hotel, and ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Our partnerships with businesses such as Japan Air, Italy Hotel, and Jamaican Restaurant, will allow you to travel and enjoy your experience at the lowest cost.  We make certain plans to hit days that have specials and specials so you can get the best deals, and get the full experience.",
    2)

# --- Step 2: the text just inserted should live in its own run (matching
# the source edit, which added it as a separate w:r). Toggling a character
# property on just that span forces the engine to break it into a distinct
# run from the text that precedes it, then clearing the property again
# keeps the split without leaving any visible formatting behind. ---
$find2 = $d.Content.Find
$found = $find2.Execute(
    "Jamaican Restaurant, will allow you to travel and enjoy your experience at the lowest cost.  We make certain plans to hit days that have specials and specials so you can get the best deals, and get the full experience.",
    $true, $false, $false, $false, $false, $true, 1, $false)
$newRun = $find2.Parent
$newRun.Font.Bold = $true
$newRun.Font.Bold = $false

# --- Step 3: split the paragraph right after that sentence (and before the
# _GoBack bookmark) into a new paragraph. ---
$newRun.Collapse(0)
$newRun.InsertParagraphAfter()

# --- Step 4: fill the new (now last) paragraph with the closing text. The
# _GoBack bookmark remains anchored at the very end of the story, i.e. at
# the end of this new paragraph. ---
$lastPara = $d.Paragraphs.Last
$closeRange = $lastPara.Range
$closeRange.End = $lastPara.Range.End - 1
$closeRange.InsertAfter("We hope that we can give you the travel experience of a lifetime.  If you have any questions feel free to contact us by going to our contact us page.  If want to book your travel now, go ahead and head over to our destinations to browse our great travel deals.  We hope we can give you the best experience of a life time.  So we have a question to ask you, where’s your escape.")
